$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.191.71"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.643.08"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'217.23"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "'19.92"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.874.59"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.661.96"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "'4.15"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'0.544"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "'67.37"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "27.193.40"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'219.06"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'147.75"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'15.76"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "1.266.11"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "'0.0178"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").Value = "'0.544"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "'0.847"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +5.84%  "
$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "1.785.04"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'61.70"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").Value = "'91.66"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "'0.0974"
$ws.Range("E51").Value = "  +0.24%  "
